$p = $ppt.ActivePresentation

# --- Update the auto-date placeholder text cache on the slide master and
# --- every slide layout (8/1/2024 -> 9/3/2024).
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "9/3/2024"
    }
}
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "9/3/2024"
        }
    }
}

# --- Slide 8: fill in the title and body placeholders with the new
# --- "read me" style notes about the Arume data presentation.
$slide8 = $p.Slides.Item(8)

$title = $slide8.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Data presentation in Arume"

$body = $slide8.Shapes.Item(2)
$btr = $body.TextFrame.TextRange
$btr.Text = "It is similar to the openiris:"
[void]$btr.InsertAfter("`rDownward is ")
$para2 = $btr.Paragraphs(2, 1)
$para2.IndentLevel = 2
$boldStart = $para2.Start + 12
[void]$btr.InsertAfter("positive")
[void]$btr.InsertAfter(" y and ")
[void]$btr.InsertAfter("to the ")
[void]$btr.InsertAfter("right is positive x")
$boldRange = $btr.Characters($boldStart, 8)
$boldRange.Font.Bold = $true
